# Swap the "Integral" / "Office Theme" colour schemes that live in this deck.
#
# Before: slide master's theme (ppt/theme/theme1.xml) = "Integral" colours,
#         notes master's theme (ppt/theme/theme2.xml)  = "Office Theme" colours.
# After:  slide master's theme should carry the "Office Theme" colours
#         (the notes-master theme part is not reachable through the
#         PowerPoint object model, so it is left untouched).
#
# PowerPoint exposes the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink -- in that fixed order) through
#   Design.SlideMaster.Theme.ThemeColorScheme.Item(n).RGB
# RGB is the usual VBA long: R + G*256 + B*65536.

function ToVbaRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Target colour values = the "Office Theme" scheme.
$tcs.Item(1).RGB  = ToVbaRGB("000000")  # dk1
$tcs.Item(2).RGB  = ToVbaRGB("FFFFFF")  # lt1
$tcs.Item(3).RGB  = ToVbaRGB("44546A")  # dk2
$tcs.Item(4).RGB  = ToVbaRGB("E7E6E6")  # lt2
$tcs.Item(5).RGB  = ToVbaRGB("5B9BD5")  # accent1
$tcs.Item(6).RGB  = ToVbaRGB("ED7D31")  # accent2
$tcs.Item(7).RGB  = ToVbaRGB("A5A5A5")  # accent3
$tcs.Item(8).RGB  = ToVbaRGB("FFC000")  # accent4
$tcs.Item(9).RGB  = ToVbaRGB("4472C4")  # accent5
$tcs.Item(10).RGB = ToVbaRGB("70AD47")  # accent6
$tcs.Item(11).RGB = ToVbaRGB("0563C1")  # hlink
$tcs.Item(12).RGB = ToVbaRGB("954F72")  # folHlink
